{"js": "// Resume content update:\n//   1. Project title \"Huskerly\" -> \"GitMarks\"\n//   2. Project description rewritten from the real-time messaging platform\n//      blurb to the GitMarks (code grading platform) blurb.\n\n// 1) Rename the project title \"Huskerly\" -> \"GitMarks\".\nconst titleResults = context.document.body.search(\"Huskerly\", { matchCase: true });\ntitleResults.load(\"text\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"GitMarks\", Word.InsertLocation.replace);\n}\n\n// 2) Replace the old project description with the new one.\nconst oldDescription =\n  \"Developed a feature-rich real-time messaging platform following a \" +\n  \"microservice architecture. Deployed using Kubernetes on AWS through a \" +\n  \"Github Actions CI/CD pipeline and Terraform for infrastructure-as-code\";\nconst newDescription =\n  \"A feature-rich code grading platform\\u2014extending the concept of GitHub \" +\n  \"Classroom to deliver an academia-to-industry submission model for \" +\n  \"students. Developed on behalf of NEU and soon-to-be piloted by the college\";\n\nconst descResults = context.document.body.search(oldDescription, { matchCase: true });\ndescResults.load(\"text\");\nawait context.sync();\n\nif (descResults.items.length > 0) {\n  descResults.items[0].insertText(newDescription, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Resume content update:\n#   1. Project title \"Huskerly\" -> \"GitMarks\"\n#   2. Project description rewritten from the real-time messaging platform\n#      blurb to the GitMarks (code grading platform) blurb.\n\n$d = $word.ActiveDocument\n\n# 1) Rename the project title \"Huskerly\" -> \"GitMarks\".\n$titleRange = $d.Content\n$titleRange.Find.MatchCase = $true\n$titleRange.Find.MatchWholeWord = $false\n$titleFound = $titleRange.Find.Execute(\"Huskerly\")\nif ($titleFound) {\n    $titleRange.Text = \"GitMarks\"\n}\n\n# 2) Replace the old project description with the new one.\n$oldDescription = \"Developed a feature-rich real-time messaging platform following a microservice architecture. Deployed using Kubernetes on AWS through a Github Actions CI/CD pipeline and Terraform for infrastructure-as-code\"\n$newDescription = \"A feature-rich code grading platform\" + [char]0x2014 + \"extending the concept of GitHub Classroom to deliver an academia-to-industry submission model for students. Developed on behalf of NEU and soon-to-be piloted by the college\"\n\n$descRange = $d.Content\n$descRange.Find.MatchCase = $true\n$descRange.Find.MatchWholeWord = $false\n$descFound = $descRange.Find.Execute($oldDescription)\nif ($descFound) {\n    $descRange.Text = $newDescription\n}\n"}
